$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Insert a new "2022-Q1" worksheet right before the "总计" sheet
#    (mirrors the per-quarter fund-holding-detail sheets such as
#    "2021-Q4") and populate it with the new fund holdings.
# ---------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$ws = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$ws.Name = "2022-Q1"

# Make sure text-like numeric values (fund codes, percentages, ...)
# are stored as text, not auto-converted to numbers.
$ws.Range("B1:G3").NumberFormat = "@"

# Copy header formatting (bold/border/centered style) from an
# existing sheet that uses the same column layout.
$src.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

# Copy the row-index column formatting too.
$src.Range("A2:A3").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "580006"
$ws.Range("C2").Value = "东吴新经济混合"
$ws.Range("D2").Value = "1.17"
$ws.Range("E2").Value = "91.56"
$ws.Range("F2").Value = "4.22"
$ws.Range("G2").Value = "0.0494"
$ws.Range("H2").Value = 9

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "002123"
$ws.Range("C3").Value = "北信瑞丰外延增长主题灵活配置混合"
$ws.Range("D3").Value = "0.17"
$ws.Range("E3").Value = "94.48"
$ws.Range("F3").Value = "4.94"
$ws.Range("G3").Value = "0.0084"
$ws.Range("H3").Value = 8

# ---------------------------------------------------------------
# 2) Add the 2022-Q1 summary row to the "总计" sheet, shifting the
#    existing rows down by one and renumbering the index column.
#    (Re-fetch the worksheet reference by name since the sheet was
#    added after our earlier lookup.)
# ---------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert(-4121)
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.06

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
